$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header swap: average_doctor / average_doctor_old (BP1 / BQ1) ---
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# --- Updated statistic values (rows 4-13) ---
$ws.Range("E4").Value = 0.478
$ws.Range("F4").Value = 0.053
$ws.Range("G4").Value = 0.23
$ws.Range("N4").Value = 0.477
$ws.Range("O4").Value = 0.064
$ws.Range("P4").Value = 0.253
$ws.Range("Q4").Value = 0.054
$ws.Range("R4").Value = 0.036
$ws.Range("S4").Value = 0.19
$ws.Range("W4").Value = 0.384
$ws.Range("X4").Value = 0.104
$ws.Range("Y4").Value = 0.322
$ws.Range("AI4").Value = 0.409
$ws.Range("AJ4").Value = 0.094
$ws.Range("AK4").Value = 0.306
$ws.Range("AU4").Value = 0.251
$ws.Range("AV4").Value = 0.024
$ws.Range("AW4").Value = 0.156
$ws.Range("BA4").Value = 2.03
$ws.Range("BB4").Value = 0.146
$ws.Range("BC4").Value = 0.382
$ws.Range("BG4").Value = 0.715
$ws.Range("BH4").Value = 0.144
$ws.Range("BI4").Value = 0.379
$ws.Range("BM4").Value = 0.744
$ws.Range("BN4").Value = 0.064
$ws.Range("BO4").Value = 0.253
$ws.Range("BP4").Value = 0.677
$ws.Range("BQ4").Value = 0.759
$ws.Range("E5").Value = 0.6
$ws.Range("F5").Value = 0.061
$ws.Range("G5").Value = 0.247
$ws.Range("N5").Value = 0.71
$ws.Range("O5").Value = 0.077
$ws.Range("P5").Value = 0.278
$ws.Range("Q5").Value = 0.035
$ws.Range("R5").Value = 0.015
$ws.Range("S5").Value = 0.122
$ws.Range("W5").Value = 0.351
$ws.Range("Y5").Value = 0.311
$ws.Range("AI5").Value = 0.409
$ws.Range("AJ5").Value = 0.09
$ws.Range("AK5").Value = 0.3
$ws.Range("AU5").Value = 0.467
$ws.Range("AV5").Value = 0.07099999999999999
$ws.Range("AW5").Value = 0.266
$ws.Range("BA5").Value = 1.294
$ws.Range("BB5").Value = 0.075
$ws.Range("BC5").Value = 0.274
$ws.Range("BG5").Value = 0.372
$ws.Range("BM5").Value = 0.526
$ws.Range("BN5").Value = 0.048
$ws.Range("BO5").Value = 0.219
$ws.Range("BP5").Value = 0.431
$ws.Range("BQ5").Value = 0.455
$ws.Range("E6").Value = 0.532
$ws.Range("N6").Value = 0.571
$ws.Range("Q6").Value = 0.042
$ws.Range("W6").Value = 0.367
$ws.Range("AI6").Value = 0.409
$ws.Range("AU6").Value = 0.327
$ws.Range("BA6").Value = 1.573
$ws.Range("BG6").Value = 0.489
$ws.Range("BM6").Value = 0.616
$ws.Range("BP6").Value = 0.524
$ws.Range("BQ6").Value = 0.5659999999999999
$ws.Range("E7").Value = 0.571
$ws.Range("N7").Value = 0.647
$ws.Range("Q7").Value = 0.038
$ws.Range("W7").Value = 0.357
$ws.Range("AI7").Value = 0.409
$ws.Range("AU7").Value = 0.398
$ws.Range("BA7").Value = 1.392
$ws.Range("BG7").Value = 0.411
$ws.Range("BM7").Value = 0.5590000000000001
$ws.Range("BP7").Value = 0.464
$ws.Range("BQ7").Value = 0.494
$ws.Range("E8").Value = 0.701
$ws.Range("F8").Value = 0.073
$ws.Range("G8").Value = 0.271
$ws.Range("N8").Value = 0.798
$ws.Range("O8").Value = 0.06
$ws.Range("P8").Value = 0.245
$ws.Range("Q8").Value = 0.039
$ws.Range("W8").Value = 0.424
$ws.Range("X8").Value = 0.116
$ws.Range("Y8").Value = 0.341
$ws.Range("AI8").Value = 0.475
$ws.Range("AJ8").Value = 0.135
$ws.Range("AK8").Value = 0.367
$ws.Range("AU8").Value = 0.413
$ws.Range("AV8").Value = 0.079
$ws.Range("AW8").Value = 0.282
$ws.Range("BA8").Value = 1.752
$ws.Range("BB8").Value = 0.11
$ws.Range("BC8").Value = 0.331
$ws.Range("BG8").Value = 0.552
$ws.Range("BH8").Value = 0.11
$ws.Range("BI8").Value = 0.332
$ws.Range("BM8").Value = 0.68
$ws.Range("BN8").Value = 0.06
$ws.Range("BO8").Value = 0.245
$ws.Range("BP8").Value = 0.584
$ws.Range("BQ8").Value = 0.621
$ws.Range("E9").Value = 0.659
$ws.Range("F9").Value = 0.225
$ws.Range("G9").Value = 0.474
$ws.Range("N9").Value = 0.732
$ws.Range("O9").Value = 0.196
$ws.Range("P9").Value = 0.443
$ws.Range("W9").Value = 0.317
$ws.Range("X9").Value = 0.217
$ws.Range("Y9").Value = 0.465
$ws.Range("AI9").Value = 0.415
$ws.Range("AJ9").Value = 0.243
$ws.Range("AK9").Value = 0.493
$ws.Range("BA9").Value = 1.707
$ws.Range("BG9").Value = 0.585
$ws.Range("BH9").Value = 0.243
$ws.Range("BI9").Value = 0.493
$ws.Range("BM9").Value = 0.659
$ws.Range("BN9").Value = 0.225
$ws.Range("BO9").Value = 0.474
$ws.Range("BP9").Value = 0.569
$ws.Range("BQ9").Value = 0.611
$ws.Range("E10").Value = 0.805
$ws.Range("F10").Value = 0.157
$ws.Range("G10").Value = 0.396
$ws.Range("N10").Value = 0.927
$ws.Range("O10").Value = 0.068
$ws.Range("P10").Value = 0.26
$ws.Range("W10").Value = 0.537
$ws.Range("AI10").Value = 0.512
$ws.Range("AJ10").Value = 0.25
$ws.Range("AK10").Value = 0.5
$ws.Range("AU10").Value = 0.415
$ws.Range("AV10").Value = 0.243
$ws.Range("AW10").Value = 0.493
$ws.Range("BA10").Value = 2.171
$ws.Range("BB10").Value = 0.217
$ws.Range("BC10").Value = 0.465
$ws.Range("BG10").Value = 0.659
$ws.Range("BH10").Value = 0.225
$ws.Range("BI10").Value = 0.474
$ws.Range("BM10").Value = 0.829
$ws.Range("BN10").Value = 0.142
$ws.Range("BO10").Value = 0.376
$ws.Range("BP10").Value = 0.724
$ws.Range("BQ10").Value = 0.759
$ws.Range("E11").Value = 0.854
$ws.Range("F11").Value = 0.125
$ws.Range("G11").Value = 0.353
$ws.Range("N11").Value = 0.927
$ws.Range("O11").Value = 0.068
$ws.Range("P11").Value = 0.26
$ws.Range("W11").Value = 0.537
$ws.Range("AI11").Value = 0.585
$ws.Range("AJ11").Value = 0.243
$ws.Range("AK11").Value = 0.493
$ws.Range("AU11").Value = 0.585
$ws.Range("AV11").Value = 0.243
$ws.Range("AW11").Value = 0.493
$ws.Range("BA11").Value = 2.171
$ws.Range("BB11").Value = 0.217
$ws.Range("BC11").Value = 0.465
$ws.Range("BG11").Value = 0.659
$ws.Range("BH11").Value = 0.225
$ws.Range("BI11").Value = 0.474
$ws.Range("BM11").Value = 0.829
$ws.Range("BN11").Value = 0.142
$ws.Range("BO11").Value = 0.376
$ws.Range("BP11").Value = 0.724
$ws.Range("BQ11").Value = 0.766
$ws.Range("E12").Value = 1.429
$ws.Range("F12").Value = 0.873
$ws.Range("G12").Value = 0.9350000000000001
$ws.Range("N12").Value = 1.263
$ws.Range("O12").Value = 0.299
$ws.Range("P12").Value = 0.547
$ws.Range("W12").Value = 1.5
$ws.Range("X12").Value = 0.432
$ws.Range("Y12").Value = 0.657
$ws.Range("AI12").Value = 1.625
$ws.Range("AJ12").Value = 1.484
$ws.Range("AK12").Value = 1.218
$ws.Range("AU12").Value = 2.846
$ws.Range("AV12").Value = 3.361
$ws.Range("AW12").Value = 1.833
$ws.Range("BA12").Value = 3.812
$ws.Range("BB12").Value = 0.459
$ws.Range("BC12").Value = 0.678
$ws.Range("BG12").Value = 1.148
$ws.Range("BH12").Value = 0.2
$ws.Range("BI12").Value = 0.448
$ws.Range("BM12").Value = 1.235
$ws.Range("BN12").Value = 0.239
$ws.Range("BO12").Value = 0.489
$ws.Range("BP12").Value = 1.271
$ws.Range("BQ12").Value = 1.254
$ws.Range("E13").Value = 1.4
$ws.Range("F13").Value = 0.294
$ws.Range("G13").Value = 0.542
$ws.Range("N13").Value = 1.723
$ws.Range("O13").Value = 0.486
$ws.Range("P13").Value = 0.697
$ws.Range("W13").Value = 0.972
$ws.Range("X13").Value = 0.192
$ws.Range("Y13").Value = 0.438
$ws.Range("AI13").Value = 1.159
$ws.Range("AJ13").Value = 0.312
$ws.Range("AK13").Value = 0.5580000000000001
$ws.Range("AU13").Value = 2.039
$ws.Range("AV13").Value = 0.339
$ws.Range("AW13").Value = 0.582
$ws.Range("BA13").Value = 2.188
$ws.Range("BB13").Value = 0.287
$ws.Range("BC13").Value = 0.536
$ws.Range("BG13").Value = 0.54
$ws.Range("BH13").Value = 0.052
$ws.Range("BI13").Value = 0.228
$ws.Range("BM13").Value = 0.794
$ws.Range("BN13").Value = 0.165
$ws.Range("BO13").Value = 0.406
$ws.Range("BP13").Value = 0.729
$ws.Range("BQ13").Value = 0.666
